$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the values of columns B, C, D for rows 2-11:
# new B = old C, new C = old D, new D = old B
for ($r = 2; $r -le 11; $r++) {
    $oldB = $ws.Cells.Item($r, 2).Value2
    $oldC = $ws.Cells.Item($r, 3).Value2
    $oldD = $ws.Cells.Item($r, 4).Value2

    $ws.Cells.Item($r, 2).Value2 = $oldC
    $ws.Cells.Item($r, 3).Value2 = $oldD
    $ws.Cells.Item($r, 4).Value2 = $oldB
}

# Update the active selection on the sheet
$ws.Range("E6").Select()
